$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to reflect the new "through" date
$ws.Name = "Through 2022-08-30"

# Update the August row label text
$ws.Range("A9").Value = "August (through 08-30)"

# Update August row (row 9) values
$ws.Range("B9").Value = 31
$ws.Range("C9").Value = 78
$ws.Range("D9").Value = 86
$ws.Range("E9").Value = 65
$ws.Range("F9").Value = 42
$ws.Range("G9").Value = 161
$ws.Range("H9").Value = 156
$ws.Range("I9").Value = 163

# Update Total row (row 10) values
$ws.Range("B10").Value = 193
$ws.Range("C10").Value = 380
$ws.Range("D10").Value = 551
$ws.Range("E10").Value = 490
$ws.Range("F10").Value = 346
$ws.Range("G10").Value = 782
$ws.Range("H10").Value = 1066
$ws.Range("I10").Value = 1134
